# Unity Gantt Chart edit script
# - Delete the "Skeleton GC Code" row (row 12) entirely, shifting later rows up
# - Update the highlighted period (H2) from 9 to 11
# - Update several actual-duration / percent-complete values
# - Move the active selection to F10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Planner")

# Delete the entire "Skeleton GC Code" row (row 12). This shifts rows 13:27
# up by one, automatically updating the dimension, conditional formatting
# ranges, and the shared-string table.
$ws.Rows(12).Delete()

# Highlighted period selector
$ws.Range("H2").Value = 11

# Row 9 (Acquire VR Setup): actual duration & percent complete
$ws.Range("F9").Value = 9
$ws.Range("G9").Value = 0.5

# Row 10 (Acquire Vision Software): actual duration & percent complete
$ws.Range("F10").Value = 8
$ws.Range("G10").Value = 0.8

# Row 11 (Skeleton VR Code): percent complete
$ws.Range("G11").Value = 1

# Row 13 (VR Interface): actual duration
$ws.Range("F13").Value = 7

# Row 14 (Vision Processing Code (GC)): actual start, actual duration, percent complete
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = 2
$ws.Range("G14").Value = 0.05

# The row delete above does not automatically shrink the conditional
# formatting sqref ranges on this runtime, so re-apply them explicitly to
# match the new (one-row-shorter) data range.
$dataCf = $ws.Range("H5:BD27").FormatConditions
$dataCf.Item(1).ModifyAppliesToRange($ws.Range("H5:BD26"))

$stripeCf = $ws.Range("B28:BD28").FormatConditions
$stripeCf.Item(1).ModifyAppliesToRange($ws.Range("B27:BD27"))

# Move selection as recorded in the workbook view
$ws.Range("F10").Select()
